$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("M2").Value = 0.259369
$ws.Range("N2").Value = 0.778107
$ws.Range("O2").Value = 0.0514155333512404
$ws.Range("P2").Value = 0.0514155333512404
$ws.Range("Q2").Value = 0.04245049194833334
$ws.Range("R2").Value = 0.382054427535
$ws.Range("S2").Value = 0.0004442041434256985
$ws.Range("T2").Value = 0.0004442041434256985
$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("O3").Value = 0.5982999525231611
$ws.Range("P3").Value = 0.5982999525231611
$ws.Range("Q3").Value = 0.4939777079383333
$ws.Range("R3").Value = 4.445799371445
$ws.Range("S3").Value = 0.005169008286010033
$ws.Range("T3").Value = 0.005169008286010034
$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 1.767033
$ws.Range("N4").Value = 5.301099
$ws.Range("O4").Value = 0.3502845141255985
$ws.Range("P4").Value = 0.3502845141255985
$ws.Range("Q4").Value = 0.289207346055
$ws.Range("R4").Value = 2.602866114495
$ws.Range("S4").Value = 0.003026280627869722
$ws.Range("T4").Value = 0.003026280627869723
$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("M5").Value = 0.259369
$ws.Range("N5").Value = 0.778107
$ws.Range("O5").Value = 0.0514155333512404
$ws.Range("P5").Value = 0.0514155333512404
$ws.Range("Q5").Value = 3.974266056794001
$ws.Range("R5").Value = 35.76839451114601
$ws.Range("S5").Value = 0.04158692558033639
$ws.Range("T5").Value = 0.04158692558033639
$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("O6").Value = 0.5982999525231611
$ws.Range("P6").Value = 0.5982999525231611
$ws.Range("Q6").Value = 46.246786488638
$ws.Range("R6").Value = 416.2210783977421
$ws.Range("S6").Value = 0.4839287658522213
$ws.Range("T6").Value = 0.4839287658522214
$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 1.767033
$ws.Range("N7").Value = 5.301099
$ws.Range("O7").Value = 0.3502845141255985
$ws.Range("P7").Value = 0.3502845141255985
$ws.Range("Q7").Value = 27.075939195258
$ws.Range("R7").Value = 243.683452757322
$ws.Range("S7").Value = 0.2833240281953454
$ws.Range("T7").Value = 0.2833240281953454
$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("M8").Value = 0.259369
$ws.Range("N8").Value = 0.778107
$ws.Range("O8").Value = 0.0514155333512404
$ws.Range("P8").Value = 0.0514155333512404
$ws.Range("Q8").Value = 0.8968231308153335
$ws.Range("R8").Value = 8.071408177338
$ws.Range("S8").Value = 0.009384403627478315
$ws.Range("T8").Value = 0.009384403627478314
$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("O9").Value = 0.5982999525231611
$ws.Range("P9").Value = 0.5982999525231611
$ws.Range("Q9").Value = 10.43593641094733
$ws.Range("R9").Value = 93.923427698526
$ws.Range("S9").Value = 0.1092021783849297
$ws.Range("T9").Value = 0.1092021783849297
$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 1.767033
$ws.Range("N10").Value = 5.301099
$ws.Range("O10").Value = 0.3502845141255985
$ws.Range("P10").Value = 0.3502845141255985
$ws.Range("Q10").Value = 6.109890030474
$ws.Range("R10").Value = 54.989010274266
$ws.Range("S10").Value = 0.06393420530238342
$ws.Range("T10").Value = 0.06393420530238342
